$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 used to hold the header labels (testcaseid / productname / result).
# The upload replaces that header row: A1 becomes a plain numeric value (3)
# and the B1/C1 header cells are removed entirely.
$ws.Range("A1").Value = 3
$ws.Range("B1").ClearContents()
$ws.Range("C1").ClearContents()

# Re-select the top-left cell (selection reverts to A1 after the header row
# is cleared out).
$ws.Range("A1").Select()
